# Weekly update: insert a new price report row for "Papa" (Terminal La
# Palmera de La Serena) above the existing data, pushing the rest of the
# table down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 323; Excel shifts rows 323:343 down to 324:344 and
# grows the used range / sheet dimension automatically.
$ws.Rows(323).Insert()

# Populate the newly inserted row with this week's figures.
$ws.Range("A323").Value = 8
$ws.Range("B323").Value = "Terminal La Palmera de La Serena"
$ws.Range("C323").Value = "Coquimbo"
$ws.Range("D323").Value = 44585
$ws.Range("E323").Value = 4
$ws.Range("F323").Value = 100114001
$ws.Range("G323").Value = "Papa"
$ws.Range("H323").Value = "Asterix"
$ws.Range("I323").Value = "1a (cosecha)"
$ws.Range("J323").Value = 2500
$ws.Range("K323").Value = 9500
$ws.Range("L323").Value = 10000
$ws.Range("M323").Value = 9750
$ws.Range("N323").Value = "`$/saco 25 kilos"
$ws.Range("O323").Value = "Provincia de Melipilla"
$ws.Range("P323").Value = 390
$ws.Range("Q323").Value = 25
$ws.Range("R323").Value = "Hortaliza"
